$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.37%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.63"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.84%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.71%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.945"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.63%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.154"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.08%"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.231"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.32%"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9281"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.46%"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1426"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.85%"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1922"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.67%"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09138"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.40%"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03517"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.20%"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09806"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.38%"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001393"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.66%"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005843"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.17%"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.937"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "6.12%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.46%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3430"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.63%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1322"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.51%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.650"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.47%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2423"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.21%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04364"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.39%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001228"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.94%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004360"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.46%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001299"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.13%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004004"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-9.97%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02040"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.58%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05054"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.64%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007388"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.25%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009812"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.16%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1363"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.45%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002129"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.14%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009382"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.11%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006358"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.60%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.12%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002727"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-18.79%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.12%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.12%"
